# Generate Report for Handoff
# Replace the old GUID-named source file ("a54144c9-bb52-4c89-bd13-114e3f915508")
# references with the new handoff run's file name
# ("ccb36728-a1a9-47cf-85bf-806625b840a0") and bump the associated timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "a54144c9-bb52-4c89-bd13-114e3f915508"
$newGuid = "ccb36728-a1a9-47cf-85bf-806625b840a0"

$oldZhHash = "cf9b556a3aea634ecd7cf3a1bf54a552571530fc"
$newZhHash = "9e010d677b1cd3dfea195a1171cd60babea09a25"

$oldDeHash = "cf9b556a3aea634ecd7cf3a1bf54a552571530fc"
$newDeHash = "9e010d677b1cd3dfea195a1171cd60babea09a25"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet (A1:G2) ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 19:01:55"

$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-22 19:01:50"

$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"

$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
